# Elimina EC anteriores y se agregan nuevos, se modifica base de datos
#
# The "Periodo Mora" data block (rows 16-24, columns B:G) is re-sorted so
# that each worker's periods are grouped together (CARLOS ENRIQUE GODOY
# RIAÑO first, then ANGELICA SUSANA GARCIA PETRO) and ordered from the
# most recent period (1903) down to the oldest, replacing the previous
# interleaved ordering. The "Valor Mora" for period 1903 is 26041 while
# all other periods keep 31249.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# TipoDoc, NumDoc, Nombre, Periodo, ValorMora, SalarioBasico
$data = @(
    @("CC", "86050699",   "CARLOS ENRIQUE GODOY RIAÑO",   "1903", 26041, 781242),
    @("CC", "86050699",   "CARLOS ENRIQUE GODOY RIAÑO",   "1902", 31249, 781242),
    @("CC", "86050699",   "CARLOS ENRIQUE GODOY RIAÑO",   "1901", 31249, 781242),
    @("CC", "86050699",   "CARLOS ENRIQUE GODOY RIAÑO",   "1811", 31249, 781242),
    @("CC", "86050699",   "CARLOS ENRIQUE GODOY RIAÑO",   "1810", 31249, 781242),
    @("CC", "1143360772", "ANGELICA SUSANA GARCIA PETRO", "1903", 26041, 781242),
    @("CC", "1143360772", "ANGELICA SUSANA GARCIA PETRO", "1902", 31249, 781242),
    @("CC", "1143360772", "ANGELICA SUSANA GARCIA PETRO", "1901", 31249, 781242),
    @("CC", "1143360772", "ANGELICA SUSANA GARCIA PETRO", "1811", 31249, 781242)
)

$startRow = 16
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rec = $data[$i]
    $ws.Range("B$row").Value = $rec[0]
    $ws.Range("C$row").Value = $rec[1]
    $ws.Range("D$row").Value = $rec[2]
    $ws.Range("E$row").Value = $rec[3]
    $ws.Range("F$row").Value = $rec[4]
    $ws.Range("G$row").Value = $rec[5]
}
